$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 6.406955
$ws.Cells.Item(2, 8).Value = 19.220865
$ws.Cells.Item(2, 9).Value = 0.2800966009992834
$ws.Cells.Item(2, 10).Value = 0.3266544289500553
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 40.91730166666667
$ws.Cells.Item(2, 14).Value = 122.751905
$ws.Cells.Item(2, 15).Value = 0.2897771170516138
$ws.Cells.Item(2, 16).Value = 0.3083463959441224
$ws.Cells.Item(2, 17).Value = 262.1553104997583
$ws.Cells.Item(2, 18).Value = 2359.397794497825
$ws.Cells.Item(2, 19).Value = 0.08116558553352851
$ws.Cells.Item(2, 20).Value = 0.100722715885935

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 6.406955
$ws.Cells.Item(3, 8).Value = 19.220865
$ws.Cells.Item(3, 9).Value = 0.2800966009992834
$ws.Cells.Item(3, 10).Value = 0.3266544289500553
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 34.738136
$ws.Cells.Item(3, 14).Value = 104.214408
$ws.Cells.Item(3, 15).Value = 0.2460161470038337
$ws.Cells.Item(3, 16).Value = 0.2617811683839066
$ws.Cells.Item(3, 17).Value = 222.56567413588
$ws.Cells.Item(3, 18).Value = 2003.09106722292
$ws.Cells.Item(3, 19).Value = 0.06890828656671387
$ws.Cells.Item(3, 20).Value = 0.0855119780683233

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 6.406955
$ws.Cells.Item(4, 8).Value = 19.220865
$ws.Cells.Item(4, 9).Value = 0.2800966009992834
$ws.Cells.Item(4, 10).Value = 0.3266544289500553
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 18.806737
$ws.Cells.Item(4, 14).Value = 56.42021099999999
$ws.Cells.Item(4, 15).Value = 0.1331896729995656
$ws.Cells.Item(4, 16).Value = 0.1417246332776418
$ws.Cells.Item(4, 17).Value = 120.493917655835
$ws.Cells.Item(4, 18).Value = 1084.445258902515
$ws.Cells.Item(4, 19).Value = 0.03730597469538437
$ws.Cells.Item(4, 20).Value = 0.04629497915146409

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 6.406955
$ws.Cells.Item(5, 8).Value = 19.220865
$ws.Cells.Item(5, 9).Value = 0.2800966009992834
$ws.Cells.Item(5, 10).Value = 0.3266544289500553
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 21.229913
$ws.Cells.Item(5, 14).Value = 63.689739
$ws.Cells.Item(5, 15).Value = 0.1503506520179033
$ws.Cells.Item(5, 16).Value = 0.1599853092240957
$ws.Cells.Item(5, 17).Value = 136.019097244915
$ws.Cells.Item(5, 18).Value = 1224.171875204235
$ws.Cells.Item(5, 19).Value = 0.04211270658824078
$ws.Cells.Item(5, 20).Value = 0.05225990982499498

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 6.406955
$ws.Cells.Item(6, 8).Value = 19.220865
$ws.Cells.Item(6, 9).Value = 0.2800966009992834
$ws.Cells.Item(6, 10).Value = 0.3266544289500553
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 13).Value = 25.510579
$ws.Cells.Item(6, 14).Value = 51.021158
$ws.Cells.Item(6, 15).Value = 0.1806664109270835
$ws.Cells.Item(6, 16).Value = 0.1281624931702333
$ws.Cells.Item(6, 17).Value = 163.445131676945
$ws.Cells.Item(6, 18).Value = 980.67079006167
$ws.Cells.Item(6, 19).Value = 0.0506040476154159
$ws.Cells.Item(6, 20).Value = 0.04186484601933792

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 6.686451000000001
$ws.Cells.Item(7, 8).Value = 20.059353
$ws.Cells.Item(7, 9).Value = 0.2923154911886005
$ws.Cells.Item(7, 10).Value = 0.3409043505233807
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 40.91730166666667
$ws.Cells.Item(7, 14).Value = 122.751905
$ws.Cells.Item(7, 15).Value = 0.2897771170516138
$ws.Cells.Item(7, 16).Value = 0.3083463959441224
$ws.Cells.Item(7, 17).Value = 273.591532646385
$ws.Cells.Item(7, 18).Value = 2462.323793817465
$ws.Cells.Item(7, 19).Value = 0.08470634030615908
$ws.Cells.Item(7, 20).Value = 0.1051166278455562

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 6.686451000000001
$ws.Cells.Item(8, 8).Value = 20.059353
$ws.Cells.Item(8, 9).Value = 0.2923154911886005
$ws.Cells.Item(8, 10).Value = 0.3409043505233807
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 34.738136
$ws.Cells.Item(8, 14).Value = 104.214408
$ws.Cells.Item(8, 15).Value = 0.2460161470038337
$ws.Cells.Item(8, 16).Value = 0.2617811683839066
$ws.Cells.Item(8, 17).Value = 232.274844195336
$ws.Cells.Item(8, 18).Value = 2090.473597758024
$ws.Cells.Item(8, 19).Value = 0.07191433085175261
$ws.Cells.Item(8, 20).Value = 0.08924233918716745

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 6.686451000000001
$ws.Cells.Item(9, 8).Value = 20.059353
$ws.Cells.Item(9, 9).Value = 0.2923154911886005
$ws.Cells.Item(9, 10).Value = 0.3409043505233807
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 18.806737
$ws.Cells.Item(9, 14).Value = 56.42021099999999
$ws.Cells.Item(9, 15).Value = 0.1331896729995656
$ws.Cells.Item(9, 16).Value = 0.1417246332776418
$ws.Cells.Item(9, 17).Value = 125.750325420387
$ws.Cells.Item(9, 18).Value = 1131.752928783483
$ws.Cells.Item(9, 19).Value = 0.03893340468411711
$ws.Cells.Item(9, 20).Value = 0.04831454406067879

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 6.686451000000001
$ws.Cells.Item(10, 8).Value = 20.059353
$ws.Cells.Item(10, 9).Value = 0.2923154911886005
$ws.Cells.Item(10, 10).Value = 0.3409043505233807
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 21.229913
$ws.Cells.Item(10, 14).Value = 63.689739
$ws.Cells.Item(10, 15).Value = 0.1503506520179033
$ws.Cells.Item(10, 16).Value = 0.1599853092240957
$ws.Cells.Item(10, 17).Value = 141.952773008763
$ws.Cells.Item(10, 18).Value = 1277.574957078867
$ws.Cells.Item(10, 19).Value = 0.04394982469513976
$ws.Cells.Item(10, 20).Value = 0.05453968793432255

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 6.686451000000001
$ws.Cells.Item(11, 8).Value = 20.059353
$ws.Cells.Item(11, 9).Value = 0.2923154911886005
$ws.Cells.Item(11, 10).Value = 0.3409043505233807
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 13).Value = 25.510579
$ws.Cells.Item(11, 14).Value = 51.021158
$ws.Cells.Item(11, 15).Value = 0.1806664109270835
$ws.Cells.Item(11, 16).Value = 0.1281624931702333
$ws.Cells.Item(11, 17).Value = 170.575236465129
$ws.Cells.Item(11, 18).Value = 1023.451418790774
$ws.Cells.Item(11, 19).Value = 0.05281159065143197
$ws.Cells.Item(11, 20).Value = 0.0436911514956556

$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 7).Value = 9.780684500000001
$ws.Cells.Item(12, 8).Value = 19.561369
$ws.Cells.Item(12, 9).Value = 0.4275879078121161
$ws.Cells.Item(12, 10).Value = 0.332441220526564
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 40.91730166666667
$ws.Cells.Item(12, 14).Value = 122.751905
$ws.Cells.Item(12, 15).Value = 0.2897771170516138
$ws.Cells.Item(12, 16).Value = 0.3083463959441224
$ws.Cells.Item(12, 17).Value = 400.1992181929909
$ws.Cells.Item(12, 18).Value = 2401.195309157945
$ws.Cells.Item(12, 19).Value = 0.1239051912119262
$ws.Cells.Item(12, 20).Value = 0.1025070522126312

$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 7).Value = 9.780684500000001
$ws.Cells.Item(13, 8).Value = 19.561369
$ws.Cells.Item(13, 9).Value = 0.4275879078121161
$ws.Cells.Item(13, 10).Value = 0.332441220526564
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 34.738136
$ws.Cells.Item(13, 14).Value = 104.214408
$ws.Cells.Item(13, 15).Value = 0.2460161470038337
$ws.Cells.Item(13, 16).Value = 0.2617811683839066
$ws.Cells.Item(13, 17).Value = 339.762748334092
$ws.Cells.Item(13, 18).Value = 2038.576490004552
$ws.Cells.Item(13, 19).Value = 0.1051935295853673
$ws.Cells.Item(13, 20).Value = 0.08702685112841589

$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 7).Value = 9.780684500000001
$ws.Cells.Item(14, 8).Value = 19.561369
$ws.Cells.Item(14, 9).Value = 0.4275879078121161
$ws.Cells.Item(14, 10).Value = 0.332441220526564
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 18.806737
$ws.Cells.Item(14, 14).Value = 56.42021099999999
$ws.Cells.Item(14, 15).Value = 0.1331896729995656
$ws.Cells.Item(14, 16).Value = 0.1417246332776418
$ws.Cells.Item(14, 17).Value = 183.9427610714765
$ws.Cells.Item(14, 18).Value = 1103.656566428859
$ws.Cells.Item(14, 19).Value = 0.05695029362006416
$ws.Cells.Item(14, 20).Value = 0.04711511006549893

$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 7).Value = 9.780684500000001
$ws.Cells.Item(15, 8).Value = 19.561369
$ws.Cells.Item(15, 9).Value = 0.4275879078121161
$ws.Cells.Item(15, 10).Value = 0.332441220526564
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 21.229913
$ws.Cells.Item(15, 14).Value = 63.689739
$ws.Cells.Item(15, 15).Value = 0.1503506520179033
$ws.Cells.Item(15, 16).Value = 0.1599853092240957
$ws.Cells.Item(15, 17).Value = 207.6430810154485
$ws.Cells.Item(15, 18).Value = 1245.858486092691
$ws.Cells.Item(15, 19).Value = 0.06428812073452281
$ws.Cells.Item(15, 20).Value = 0.05318571146477812

$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 7).Value = 9.780684500000001
$ws.Cells.Item(16, 8).Value = 19.561369
$ws.Cells.Item(16, 9).Value = 0.4275879078121161
$ws.Cells.Item(16, 10).Value = 0.332441220526564
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 13).Value = 25.510579
$ws.Cells.Item(16, 14).Value = 51.021158
$ws.Cells.Item(16, 15).Value = 0.1806664109270835
$ws.Cells.Item(16, 16).Value = 0.1281624931702333
$ws.Cells.Item(16, 17).Value = 249.5109246113255
$ws.Cells.Item(16, 18).Value = 998.0436984453021
$ws.Cells.Item(16, 19).Value = 0.07725077266023568
$ws.Cells.Item(16, 20).Value = 0.04260649565523979
